$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.487.48"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "3.979.80"
$ws.Range("E3").Value = "  -2.34%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.38"
$ws.Range("E5").Value = "  +7.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.19"
$ws.Range("E6").Value = "  +11.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.683"
$ws.Range("E7").Value = "  -1.93%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.787"
$ws.Range("E9").Value = "  +1.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  +6.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.14"
$ws.Range("E11").Value = "  +3.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000336"
$ws.Range("E12").Value = "  +2.09%  "

$ws.Range("E13").Value = "  +1.27%  "

$ws.Range("D14").Value = "4.617.20"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "3.976.95"
$ws.Range("E15").Value = "  -2.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("E16").Value = "  -2.33%  "

$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.61"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").Value = "73.392.55"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "456.40"
$ws.Range("E21").Value = "  +1.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  +8.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.29"
$ws.Range("E23").Value = "  -2.24%  "

$ws.Range("E24").Value = "  -4.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.26"
$ws.Range("E25").Value = "  -3.68%  "

$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.09"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.97"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.54"
$ws.Range("E29").Value = "  -4.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.37"
$ws.Range("E30").Value = "  -2.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.92"
$ws.Range("E32").Value = "  +2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000105"
$ws.Range("E33").Value = "  +14.34%  "

$ws.Range("E34").Value = "  -3.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.14"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.77"
$ws.Range("E36").Value = "  +3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "648.94"
$ws.Range("E37").Value = "  -5.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.431"
$ws.Range("E38").Value = "  -4.16%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0484"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.19"
$ws.Range("E44").Value = "  +37.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.53"
$ws.Range("E45").Value = "  -6.35%  "

$ws.Range("E46").Value = "  -6.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.150"
$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000299"
$ws.Range("E48").Value = "  +6.90%  "

$ws.Range("E49").Value = "  +2.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.57"
$ws.Range("E50").Value = "  -4.79%  "

$ws.Range("D51").Value = "2.822.46"
$ws.Range("E51").Value = "  +1.32%  "
